$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.45349999999999
$ws.Range("C7").Value = -12.89119999999999
$ws.Range("A8").Value = -22.32050000000002
$ws.Range("A10").Value = -21.44979999999999
$ws.Range("A12").Value = -21.57810000000001
$ws.Range("C15").Value = -14.17929999999998
$ws.Range("A18").Value = -21.4515
$ws.Range("C18").Value = -10.3462
$ws.Range("E18").Value = 18.77640000000002
$ws.Range("E19").Value = 16.49410000000002
$ws.Range("C20").Value = -11.90920000000001
$ws.Range("E27").Value = 16.4761
$ws.Range("C29").Value = -11.4689
$ws.Range("C30").Value = -12.5365
$ws.Range("C31").Value = -12.7168
$ws.Range("E31").Value = 16.4785
$ws.Range("A37").Value = -20.12420000000001
$ws.Range("E38").Value = 16.17419999999998
$ws.Range("C40").Value = -12.41050000000001
$ws.Range("E42").Value = 16.2723
$ws.Range("E44").Value = 16.37339999999999
$ws.Range("E47").Value = 16.44909999999999
$ws.Range("C50").Value = -13.5201
$ws.Range("A55").Value = -22.3486
$ws.Range("E58").Value = 16.18530000000002
$ws.Range("E65").Value = 17.25880000000001
$ws.Range("A68").Value = -21.6668
$ws.Range("C68").Value = -11.3927
$ws.Range("E73").Value = 17.3797
$ws.Range("C76").Value = -12.18530000000001
$ws.Range("A77").Value = -20.81859999999999
$ws.Range("A78").Value = -20.47179999999998
$ws.Range("A81").Value = -21.7619
$ws.Range("A82").Value = -21.99490000000001
$ws.Range("C87").Value = -13.73559999999998
$ws.Range("C88").Value = -12.96829999999999
$ws.Range("E90").Value = 16.3978
$ws.Range("E94").Value = 18.86210000000002
$ws.Range("E95").Value = 18.21050000000001
$ws.Range("C96").Value = -13.0194
$ws.Range("C98").Value = -12.1268
$ws.Range("C101").Value = -12.1743
$ws.Range("E101").Value = 16.62450000000001
$ws.Range("C102").Value = -13.0364
